$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------------------------
# Approach note: this runtime's Rows().Insert()/Columns().Insert() shift cells correctly
# but leave a stray unused cellXf behind (a format inherited from the row/col above/left
# that ends up unreferenced once we overwrite the new band with real content). To keep
# styles.xml pristine we instead shift content manually with Copy + PasteSpecial. Also:
# PasteSpecial(xlPasteValues) from an EMPTY source cell is a no-op on the destination
# (it does not clear it), so every values-only paste below is preceded by ClearContents
# on the destination range.
# ------------------------------------------------------------------------------------

# --- 1. Shift column D ("Date of Completion") to E, opening up a blank column D ---
$ws.Range("D1:D9").Copy()
[void]$ws.Range("E1:E9").PasteSpecial(-4122)
[void]$excel.CutCopyMode
$ws.Range("E1:E9").ClearContents()
$ws.Range("D1:D9").Copy()
[void]$ws.Range("E1:E9").PasteSpecial(-4163)
[void]$excel.CutCopyMode

# --- 2. Shift rows 2-9 down to rows 3-10 (bottom-up so sources aren't clobbered) ---
$ws.Range("A9:E9").Copy()
[void]$ws.Range("A10:E10").PasteSpecial(-4122)
[void]$excel.CutCopyMode
$ws.Range("A10:E10").ClearContents()
$ws.Range("A9:E9").Copy()
[void]$ws.Range("A10:E10").PasteSpecial(-4163)
[void]$excel.CutCopyMode

$ws.Range("A8:E8").Copy()
[void]$ws.Range("A9:E9").PasteSpecial(-4122)
[void]$excel.CutCopyMode
$ws.Range("A9:E9").ClearContents()
$ws.Range("A8:E8").Copy()
[void]$ws.Range("A9:E9").PasteSpecial(-4163)
[void]$excel.CutCopyMode

$ws.Range("A7:E7").Copy()
[void]$ws.Range("A8:E8").PasteSpecial(-4122)
[void]$excel.CutCopyMode
$ws.Range("A8:E8").ClearContents()
$ws.Range("A7:E7").Copy()
[void]$ws.Range("A8:E8").PasteSpecial(-4163)
[void]$excel.CutCopyMode

$ws.Range("A6:E6").Copy()
[void]$ws.Range("A7:E7").PasteSpecial(-4122)
[void]$excel.CutCopyMode
$ws.Range("A7:E7").ClearContents()
$ws.Range("A6:E6").Copy()
[void]$ws.Range("A7:E7").PasteSpecial(-4163)
[void]$excel.CutCopyMode

$ws.Range("A5:E5").Copy()
[void]$ws.Range("A6:E6").PasteSpecial(-4122)
[void]$excel.CutCopyMode
$ws.Range("A6:E6").ClearContents()
$ws.Range("A5:E5").Copy()
[void]$ws.Range("A6:E6").PasteSpecial(-4163)
[void]$excel.CutCopyMode

$ws.Range("A4:E4").Copy()
[void]$ws.Range("A5:E5").PasteSpecial(-4122)
[void]$excel.CutCopyMode
$ws.Range("A5:E5").ClearContents()
$ws.Range("A4:E4").Copy()
[void]$ws.Range("A5:E5").PasteSpecial(-4163)
[void]$excel.CutCopyMode

$ws.Range("A3:E3").Copy()
[void]$ws.Range("A4:E4").PasteSpecial(-4122)
[void]$excel.CutCopyMode
$ws.Range("A4:E4").ClearContents()
$ws.Range("A3:E3").Copy()
[void]$ws.Range("A4:E4").PasteSpecial(-4163)
[void]$excel.CutCopyMode

$ws.Range("A2:E2").Copy()
[void]$ws.Range("A3:E3").PasteSpecial(-4122)
[void]$excel.CutCopyMode
$ws.Range("A3:E3").ClearContents()
$ws.Range("A2:E2").Copy()
[void]$ws.Range("A3:E3").PasteSpecial(-4163)
[void]$excel.CutCopyMode

# --- 3. Header for new column D: "Date of Start" ---
$ws.Range("D1").Value = "Date of Start"

# --- 4. New row 2: "Setup Jetson" task. Copy formatting from row 3 (same assignee group) ---
$ws.Range("A3:E3").Copy()
[void]$ws.Range("A2:E2").PasteSpecial(-4122)
[void]$excel.CutCopyMode
$ws.Range("A2:E2").ClearContents()
$ws.Range("A2").Value = "Setup Jetson"
$ws.Range("B2").Value = "Tejas M K"
$ws.Range("C2").Value = 45999
$ws.Range("D2").Value = 46010

# --- 5. Fill "Date of Start" (col D) values for each remaining task row ---
$ws.Range("D3").Value = 46000
# Row 4 (Aruco Marker Placement & Sizing) stays blank
$ws.Range("D5").Value = 46007
$ws.Range("D6").Value = 46000
# Row 7 (Fabrication of Elevated Structure) stays blank
$ws.Range("D8").Value = 46000
# Row 9 (Map Switching Logic using Travelled Distance) stays blank
# Row 10 (Integration of Aruco Based Localization...) stays blank

# --- 6. Rename task text (dropped "via ZED Mini") ---
$ws.Range("A9").Value = "Map Switching Logic using Travelled Distance"

# --- 7. Fix a style inconsistency inherited from source data: E5 (old D4) should use the
#          numeric-date style (s=3) like its row siblings, not the plain text style (s=2). ---
$ws.Range("D5").Copy()
[void]$ws.Range("E5").PasteSpecial(-4122)
[void]$excel.CutCopyMode

# --- 8. Column D width, matching column C ---
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth

# --- 9. Selection moves to A3 ---
[void]$ws.Range("A3").Select()
